$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 12.15279999999999
$ws.Range("E7").Value = 12.0676
$ws.Range("C8").Value = -10.92619999999999
$ws.Range("A12").Value = -21.98400000000001
$ws.Range("C12").Value = -13.30560000000001
$ws.Range("C14").Value = -11.81389999999999
$ws.Range("E19").Value = 13.20349999999999
$ws.Range("E21").Value = 12.69449999999999
$ws.Range("C22").Value = -10.84949999999999
$ws.Range("E24").Value = 12.75489999999999
